$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2980
$ws.Range("I74").Value = 2950
$ws.Range("J74").Value = 3100
$ws.Range("K74").Value = 2950
$ws.Range("L74").Value = 3100
$ws.Range("M74").Value = -2014
$ws.Range("N74").Value = -4972

$ws.Range("H77").Value = 2980
$ws.Range("I77").Value = 2950
$ws.Range("J77").Value = 3100
$ws.Range("K77").Value = 14750
$ws.Range("L77").Value = 15500
$ws.Range("M77").Value = -10070
$ws.Range("N77").Value = -24860

$ws.Range("H111").Value = 4173266
$ws.Range("I111").Value = 7522.5264
$ws.Range("J111").Value = 20003092
$ws.Range("K111").Value = 22567.5792
$ws.Range("L111").Value = 60009276
$ws.Range("M111").Value = -19500.5792
$ws.Range("N111").Value = -60015410

$ws.Range("H137").Value = 1076.7097
$ws.Range("I137").Value = 1070.8125
$ws.Range("J137").Value = 1096.9286
$ws.Range("K137").Value = 3212.4375
$ws.Range("L137").Value = 3290.7858
$ws.Range("M137").Value = -662.4375
$ws.Range("N137").Value = -8390.7858

$ws.Range("H141").Value = 1210.5443
$ws.Range("I141").Value = 1071.1733
$ws.Range("J141").Value = 3823.75
$ws.Range("K141").Value = 3213.5199
$ws.Range("L141").Value = 11471.25
$ws.Range("M141").Value = 1966.4801
$ws.Range("N141").Value = -21831.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 64056.688
$ws.Range("I2").Value = 1998.7142
$ws.Range("J2").Value = 112324
$ws.Range("K2").Value = 1998.7142
$ws.Range("L2").Value = 112324
$ws.Range("M2").Value = -1885.7142
$ws.Range("N2").Value = -112550

$ws.Range("H32").Value = 2946.25
$ws.Range("I32").Value = 2680.4736
$ws.Range("J32").Value = 7996
$ws.Range("K32").Value = 2680.4736
$ws.Range("L32").Value = 7996
$ws.Range("M32").Value = -2393.4736
$ws.Range("N32").Value = -8570

$ws.Range("H45").Value = 64704.188
$ws.Range("I45").Value = 78639
$ws.Range("J45").Value = 4320
$ws.Range("K45").Value = 78639
$ws.Range("L45").Value = 4320
$ws.Range("M45").Value = -78262
$ws.Range("N45").Value = -5074

$ws.Range("H61").Value = 1031.1923
$ws.Range("I61").Value = 957.95557
$ws.Range("J61").Value = 1502
$ws.Range("K61").Value = 957.95557
$ws.Range("L61").Value = 1502
$ws.Range("M61").Value = -745.95557
$ws.Range("N61").Value = -1926

$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H116").Value = 64056.688
$ws.Range("I116").Value = 1998.7142
$ws.Range("J116").Value = 112324
$ws.Range("K116").Value = 1998.7142
$ws.Range("L116").Value = 112324
$ws.Range("M116").Value = 295.2858000000001
$ws.Range("N116").Value = -116912

$ws.Range("H132").Value = 3587.1226
$ws.Range("I132").Value = 3703.6924
$ws.Range("J132").Value = 3132.5
$ws.Range("K132").Value = 11111.0772
$ws.Range("L132").Value = 9397.5
$ws.Range("M132").Value = -8581.0772
$ws.Range("N132").Value = -14457.5

$ws.Range("H136").Value = 1031.1923
$ws.Range("I136").Value = 957.95557
$ws.Range("J136").Value = 1502
$ws.Range("K136").Value = 2873.86671
$ws.Range("L136").Value = 4506
$ws.Range("M136").Value = -323.8667099999998
$ws.Range("N136").Value = -9606

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 64056.688
$ws.Range("I3").Value = 1998.7142
$ws.Range("J3").Value = 112324
$ws.Range("K3").Value = 1998.7142
$ws.Range("L3").Value = 112324
$ws.Range("M3").Value = -1884.7142
$ws.Range("N3").Value = -112552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1175.0317
$ws.Range("I58").Value = 892.9039
$ws.Range("J58").Value = 2508.7273
$ws.Range("K58").Value = 892.9039
$ws.Range("L58").Value = 2508.7273
$ws.Range("M58").Value = -689.9039
$ws.Range("N58").Value = -2914.7273

$ws.Range("H132").Value = 3479.561
$ws.Range("I132").Value = 3529.577
$ws.Range("J132").Value = 3392.8667
$ws.Range("K132").Value = 10588.731
$ws.Range("L132").Value = 10178.6001
$ws.Range("M132").Value = -8058.731
$ws.Range("N132").Value = -15238.6001

$ws.Range("H134").Value = 1415.9131
$ws.Range("I134").Value = 1227.7333
$ws.Range("J134").Value = 1768.75
$ws.Range("K134").Value = 3683.199900000001
$ws.Range("L134").Value = 5306.25
$ws.Range("M134").Value = -1148.199900000001
$ws.Range("N134").Value = -10376.25

$ws.Range("H136").Value = 1175.0317
$ws.Range("I136").Value = 892.9039
$ws.Range("J136").Value = 2508.7273
$ws.Range("K136").Value = 2678.7117
$ws.Range("L136").Value = 7526.1819
$ws.Range("M136").Value = -128.7116999999998
$ws.Range("N136").Value = -12626.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 54.526318
$ws.Range("I38").Value = 51.57143
$ws.Range("J38").Value = 56.25
$ws.Range("K38").Value = 154.71429
$ws.Range("L38").Value = 168.75
$ws.Range("M38").Value = 192.28571
$ws.Range("N38").Value = -862.75

$ws.Range("H131").Value = 7618.383
$ws.Range("I131").Value = 742.1111
$ws.Range("J131").Value = 8477.916999999999
$ws.Range("K131").Value = 2226.3333
$ws.Range("L131").Value = 25433.751
$ws.Range("M131").Value = 2813.6667
$ws.Range("N131").Value = -35513.751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1755.3
$ws.Range("I113").Value = 2337
$ws.Range("J113").Value = 1506
$ws.Range("K113").Value = 2337
$ws.Range("L113").Value = 1506
$ws.Range("M113").Value = -167
$ws.Range("N113").Value = -5846

$ws.Range("H126").Value = 3032.75
$ws.Range("I126").Value = 3509.8
$ws.Range("J126").Value = 2237.6667
$ws.Range("K126").Value = 10529.4
$ws.Range("L126").Value = 6713.000100000001
$ws.Range("M126").Value = -8059.400000000001
$ws.Range("N126").Value = -11653.0001

$ws.Range("H132").Value = 2388.3125
$ws.Range("I132").Value = 2116.0454
$ws.Range("J132").Value = 2987.3
$ws.Range("K132").Value = 6348.1362
$ws.Range("L132").Value = 8961.900000000001
$ws.Range("M132").Value = -3818.1362
$ws.Range("N132").Value = -14021.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2238.3572
$ws.Range("I132").Value = 2160.6382
$ws.Range("J132").Value = 2644.2222
$ws.Range("K132").Value = 6481.9146
$ws.Range("L132").Value = 7932.6666
$ws.Range("M132").Value = -3951.9146
$ws.Range("N132").Value = -12992.6666

$ws.Range("H136").Value = 1120.902
$ws.Range("I136").Value = 928.1556
$ws.Range("J136").Value = 2566.5
$ws.Range("K136").Value = 2784.4668
$ws.Range("L136").Value = 7699.5
$ws.Range("M136").Value = -234.4668000000001
$ws.Range("N136").Value = -12799.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 41104.4
$ws.Range("I29").Value = 21833.334
$ws.Range("J29").Value = 70011
$ws.Range("K29").Value = 21833.334
$ws.Range("L29").Value = 70011
$ws.Range("M29").Value = -21543.334
$ws.Range("N29").Value = -70591

$ws.Range("H126").Value = 1448.2858
$ws.Range("I126").Value = 1401.2667
$ws.Range("J126").Value = 1565.8334
$ws.Range("K126").Value = 4203.800099999999
$ws.Range("L126").Value = 4697.5002
$ws.Range("M126").Value = -1733.800099999999
$ws.Range("N126").Value = -9637.5002

$ws.Range("H132").Value = 1787.4286
$ws.Range("I132").Value = 1876.9246
$ws.Range("J132").Value = 1508.4117
$ws.Range("K132").Value = 5630.7738
$ws.Range("L132").Value = 4525.2351
$ws.Range("M132").Value = -3100.7738
$ws.Range("N132").Value = -9585.2351
